$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as plain text so that numeric-looking
# strings (e.g. "1.000", "0.9996") are preserved verbatim rather than
# being coerced into numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.568.34'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.923.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.07%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4842'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.21%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2902'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06801'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.95%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '112.82'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.42'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.918.91'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.24%  '

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.486'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.41%  '

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07572'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.04%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6751'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.51%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '293.94'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.64%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.571.77'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007672'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.525'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.168.10'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9993'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.462'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.493'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.31%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.29'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.26%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.35%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.107'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.40%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1070'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.439'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.28%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.145'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.062'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04955'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7366'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.141'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.69%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02033'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.07%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.70%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.025'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.24%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '109.72'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4452'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8684'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.67%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.851'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.90%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.36'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.268'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.282'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.42%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1232'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.89%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.97'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.68%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2508'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.04%  '

Write-Output "Applied cryptos list update"